# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.318.65"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.986.78"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.36"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.36"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.15"
$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.365"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.501.17"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.61"
$ws.Range("E14").Value = "  -2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000162"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.267.25"
$ws.Range("E16").Value = "  -1.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.994.90"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.97"
$ws.Range("E18").Value = "  -1.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.17"
$ws.Range("E21").Value = "  +3.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.492"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.58"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.108.79"
$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  +0.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  +3.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.34"
$ws.Range("E29").Value = "  -3.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.89"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.20"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.47"
$ws.Range("E34").Value = "  -2.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.98"
$ws.Range("E37").Value = "  +6.53%  "

$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0658"
$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.028.20"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.89"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.652"
$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.182.42"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("E46").Value = "  -3.33%  "

$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.919"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.45"
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0848"
$ws.Range("E51").Value = "  -3.02%  "
